$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.575.82'
$ws.Range('E2').Value = '  -1.01%  '

$ws.Range('D3').Value = '2.318.09'
$ws.Range('E3').Value = '  -0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '513.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.55%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.59%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.534'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.90%  '

$ws.Range('E9').Value = '  -3.36%  '

$ws.Range('E10').Value = '  -0.18%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.26'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.336'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.51'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.94%  '

$ws.Range('D14').Value = '2.731.06'
$ws.Range('E14').Value = '  -0.44%  '

$ws.Range('D15').Value = '56.525.28'
$ws.Range('E15').Value = '  -0.80%  '

$ws.Range('E16').Value = '  -1.69%  '

$ws.Range('D17').Value = '2.309.86'
$ws.Range('E17').Value = '  -0.76%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.31%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '328.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.92%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.20%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.26%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.28'
$ws.Range('D23').Style = 'Normal'

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.81%  '

$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.164'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.35%  '

$ws.Range('E27').Value = '  +0.91%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.23%  '

$ws.Range('E29').Value = '  -2.81%  '

$ws.Range('D30').Value = '0.0₃0718'
$ws.Range('E30').Value = '  -3.69%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.86%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.30'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.37%  '

$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.93'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.64%  '

$ws.Range('E37').Value = '  -4.87%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.74%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '149.61'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.98%  '

$ws.Range('E41').Value = '  -1.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.89%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '276.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.42%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.92%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0928'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0495'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.39%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.553'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.03%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.30'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.55%  '

$ws.Range('E49').Value = '  -1.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.377'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.49%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.06'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.05%  '
